$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (row 1) to the cleaned machine-readable names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Normalize municipality/state name casing: lowercase connectors
#     (de/del/la/las/el/los/y) -> title case, matching the cleaning script ---
$ws.Range("B10").Value = "Mazapa De Madero"
$ws.Range("B12").Value = "San Cristóbal De Las Casas"
$ws.Range("B16").Value = "Hidalgo Del Parral"
$ws.Range("A22").Value = "Ciudad De México"
$ws.Range("A34").Value = "Estado De México"
$ws.Range("B34").Value = "Coacalco De Berriozábal"
$ws.Range("B35").Value = "Ecatepec De Morelos"
$ws.Range("B39").Value = "San Felipe Del Progreso"
$ws.Range("B42").Value = "Tlalnepantla De Baz"
$ws.Range("B47").Value = "Apaseo El Alto"
$ws.Range("B48").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B56").Value = "Silao De La Victoria"
$ws.Range("B58").Value = "Acapulco De Juárez"
$ws.Range("B59").Value = "Ajuchitlán Del Progreso"
$ws.Range("B60").Value = "Ayutla De Los Libres"
$ws.Range("B62").Value = "Cutzamala De Pinzón"
$ws.Range("B65").Value = "Técpan De Galeana"
$ws.Range("B74").Value = "Santiago De Anaya"
$ws.Range("B75").Value = "Tenango De Doria"
$ws.Range("B78").Value = "Ahualulco De Mercado"
$ws.Range("B83").Value = "Ojuelos De Jalisco"
$ws.Range("B98").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B112").Value = "Mier Y Noriega"
$ws.Range("B116").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B117").Value = "Oaxaca De Juárez"
$ws.Range("B118").Value = "Putla Villa De Guerrero"
$ws.Range("B143").Value = "San Salvador El Verde"
$ws.Range("B145").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B147").Value = "Landa De Matamoros"
$ws.Range("B148").Value = "Pinal De Amoles"
$ws.Range("B149").Value = "San Juan Del Río"
$ws.Range("B157").Value = "San Ciro De Acosta"
$ws.Range("B170").Value = "Soto La Marina"
$ws.Range("B183").Value = "Hueyapan De Ocampo"
$ws.Range("B186").Value = "Martínez De La Torre"
$ws.Range("B196").Value = "Zozocolco De Hidalgo"

# --- Remove obsolete footer rows (208-212): sample size / source / author notes ---
$ws.Rows("208:212").Delete()
